$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "A little stressful"
$ws.Range("D3").Value = "A little stressful"
$ws.Range("D4").Value = "Not stressful"
$ws.Range("D5").Value = "Not stressful"
$ws.Range("D6").Value = "A little stressful"
$ws.Range("D7").Value = "Not stressful"

$ws.Range("D8").Select()
